$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for a new "Turma" column between "Disciplina" (B) and "Nota" (C):
# shift the existing Nota/Preferencia columns one column to the right,
# using Copy so the original value types/number-formats/styles survive.
$ws.Range("D1:D3").Copy($ws.Range("E1:E3"))
$ws.Range("C1:C3").Copy($ws.Range("D1:D3"))

# Fill in the new "Turma" column
$ws.Range("C1").Value = "Turma"
$ws.Range("C2").Value = 1
$ws.Range("C3").Value = 2

# Match the bold header style used by the rest of row 1
$ws.Range("C1").Font.Bold = $true

# A new styled (underlined) empty cell shows up at E13, mirroring the one at G5
$ws.Range("E13").Font.Underline = $true
$ws.Range("E13").Select()
